$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("标签")

# Add new tags under the "深度学习" (column B) and "机器学习" (column C) columns
# Values are entered in this order so the shared-string table is built in the
# same sequence as the target workbook (Ensemble, Bagging, IMDB, SaveModel, TensorFlow)
$ws.Range("C9").Value = "Ensemble"
$ws.Range("C10").Value = "Bagging"
$ws.Range("B9").Value = "IMDB"
$ws.Range("B10").Value = "SaveModel"
$ws.Range("B11").Value = "TensorFlow"

# Update selection to match the author's last edited cell
$ws.Activate()
$ws.Range("B11").Select()
